$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

function Set-RowValues {
    param($RangeAddress, $Values)
    $arr = New-Object 'object[,]' 1,$Values.Length
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $ws.Range($RangeAddress).Value = $arr
}

# Row 3: quicksort_many_equal
Set-RowValues "B3:J3" @(6900, 13720, 28380, 69960, 154600, 432150, 1347540, 4532630, 16446350)

# Row 5: heapsort_many_equal
Set-RowValues "B5:J5" @(8160, 18020, 39220, 84860, 180460, 378480, 757480, 1508130, 3012890)

# Row 7: insertion_sort_many_equal
Set-RowValues "B7:J7" @(75050, 291030, 1165580, 4668220, 18541010, 74027310, 294371060, 1183581160, 4725895550)

# Update the active selection on the sheet to I9
$ws.Range("I9").Select()
